$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# to reflect the latest values from the automated data refresh.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.011.71"
$ws.Range("E2").Value = "  -0.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.411.68"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "408.47"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.76"
$ws.Range("E6").Value = "  -3.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.638"
$ws.Range("E7").Value = "  +7.80%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.732"
$ws.Range("E9").Value = "  +6.64%  "

$ws.Range("E10").Value = "  +17.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.45"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000218"
$ws.Range("E12").Value = "  +67.33%  "

$ws.Range("E13").Value = "  -0.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.961.69"
$ws.Range("E14").Value = "  +0.24%  "

$ws.Range("E15").Value = "  +5.93%  "

$ws.Range("E16").Value = "  +4.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.400.29"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("E18").Value = "  +10.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.06"
$ws.Range("E19").Value = "  +4.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "61.947.44"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "425.52"
$ws.Range("E21").Value = "  +35.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "89.23"
$ws.Range("E22").Value = "  +5.75%  "

$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.01"
$ws.Range("E24").Value = "  +0.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.24"
$ws.Range("E25").Value = "  +2.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "33.31"
$ws.Range("E26").Value = "  +12.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.81"
$ws.Range("E27").Value = "  +7.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.78"
$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.57"
$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.69"
$ws.Range("E30").Value = "  -5.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.91"
$ws.Range("E31").Value = "  +4.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.170"
$ws.Range("E32").Value = "  -3.29%  "

$ws.Range("E33").Value = "  -0.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.66"
$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("E35").Value = "  +0.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0499"
$ws.Range("E36").Value = "  +3.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.17"
$ws.Range("E37").Value = "  +4.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.35"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("E40").Value = "  +7.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("E41").Value = "  -1.41%  "

$ws.Range("E42").Value = "  +3.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.03"
$ws.Range("E43").Value = "  +2.57%  "

$ws.Range("E45").Value = "  +2.16%  "

$ws.Range("E46").Value = "  +8.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.59"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.03"
$ws.Range("E48").Value = "  +2.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.759.92"
$ws.Range("E49").Value = "  +0.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.112.35"
$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("E51").Value = "  +2.69%  "
